# Applies odds updates to rows 3, 6, and 8 of the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("G3").Value = 3.3
$ws.Range("I3").Value = 2.3
$ws.Range("L3").Value = 3
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 15
$ws.Range("AH3").Value = 8
$ws.Range("AO3").Value = 17

# --- Row 6 ---
$ws.Range("G6").Value = 1.53
$ws.Range("H6").Value = 3.9
$ws.Range("I6").Value = 6.25
$ws.Range("J6").Value = 2.1
$ws.Range("K6").Value = 2.38
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 12
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("W6").Value = 7.5
$ws.Range("Z6").Value = 11
$ws.Range("AA6").Value = 12
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 12
$ws.Range("AD6").Value = 7.5
$ws.Range("AH6").Value = 17
$ws.Range("AJ6").Value = 19
$ws.Range("AK6").Value = 67
$ws.Range("AO6").Value = 7.5
$ws.Range("AQ6").Value = 23
$ws.Range("AR6").Value = 41
$ws.Range("AT6").Value = 3.25
$ws.Range("BB6").Value = 201

# --- Row 8 ---
$ws.Range("G8").Value = 3.6
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 2.05
$ws.Range("J8").Value = 4
$ws.Range("L8").Value = 2.63
$ws.Range("W8").Value = 12
$ws.Range("X8").Value = 19
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 41
$ws.Range("AA8").Value = 26
$ws.Range("AB8").Value = 34
$ws.Range("AH8").Value = 8.5
$ws.Range("AI8").Value = 10
$ws.Range("AK8").Value = 19
$ws.Range("AL8").Value = 15
$ws.Range("AO8").Value = 19
$ws.Range("AP8").Value = 26
$ws.Range("AR8").Value = 81
$ws.Range("AX8").Value = 11
$ws.Range("AY8").Value = 19
$ws.Range("AZ8").Value = 34
